# Add a new "Player Info" worksheet before "ODI Batting" and populate it,
# then convert the MATCH_CARD_LINK url columns on the batting/bowling
# sheets into plain MATCH_CODE numbers (stored as text).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet as the first tab -----------------
$battingWsForInsert = $wb.Worksheets.Item("ODI Batting")
$playerWs = $wb.Worksheets.Add($battingWsForInsert)
$playerWs.Name = "Player Info"

# NOTE: sheet handles in this host resolve by tab position, so re-fetch the
# other sheets by name now that the tab order has shifted.
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# Header row
$playerWs.Range("A1").Value = "ID"
$playerWs.Range("B1").Value = "NAME"
$playerWs.Range("C1").Value = "BATTING_HAND"
$playerWs.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header look used on the other sheets
$headerRange = $playerWs.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data row
$playerWs.Range("B2").Value = "Nathan Ellis"
$playerWs.Range("C2").Value = "Right Handed"
$playerWs.Range("D2").Value = "Right Arm Fast Medium"

# Keep the numeric-looking ID as text (matches the scraped-data formatting)
$playerWs.Range("A2").Formula = '=TEXT(6082,"0")'
$playerWs.Range("A2").Copy() | Out-Null
$playerWs.Range("A2").PasteSpecial(-4163) | Out-Null

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ---------------------------
$battingWs.Range("D1").Value = "MATCH_CODE"

$battingCodes = @(
    @{ Row = 2; Code = "4564" },
    @{ Row = 3; Code = "4565" },
    @{ Row = 4; Code = "4567" },
    @{ Row = 5; Code = "4728" }
)
foreach ($item in $battingCodes) {
    $cell = $battingWs.Range("D" + $item.Row)
    $cell.Formula = '=TEXT(' + $item.Code + ',"0")'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE ---------------------------
$bowlingWs.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @(
    @{ Row = 2; Code = "4564" },
    @{ Row = 3; Code = "4565" },
    @{ Row = 4; Code = "4567" },
    @{ Row = 5; Code = "4728" }
)
foreach ($item in $bowlingCodes) {
    $cell = $bowlingWs.Range("B" + $item.Row)
    $cell.Formula = '=TEXT(' + $item.Code + ',"0")'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$playerWs.Range("A1").Select() | Out-Null
